# Simulates a second run of the Selenium "write web table into Excel" test
# suite (now parallelised, per the commit message) against the same
# workbook: each sheet's table is rewritten with a fresh "FineshedAt"
# timestamp, and a brand-new duplicate copy of the same table (also with
# its own fresh timestamp) is appended immediately to the right of the
# first one, since the parallel data-provider runs no longer stomp on one
# another's columns.

$wb = $excel.ActiveWorkbook

$HEADER_COLOR_INDEX = 48   # indexed palette 55 -> the same fill used by the
                           # header row of every table already in the file

# Calibrated ColumnWidth (Excel "characters" units) so the re-saved <col>
# widths line up with the four distinct column widths already used
# elsewhere in the workbook (28.36328125 / 16.1328125 / 9.23828125 /
# 11.1328125 "pixel-based" units).
$WIDE_COL_WIDTH   = 27.5     # ~28.36328125
$MEDIUM_COL_WIDTH = 15.3125  # ~16.1328125
$NARROW_COL_WIDTH = 8.28125  # ~9.23828125
$DATE_COL_WIDTH   = 10.3125  # ~11.1328125

function Style-Header($rng) {
    $rng.Interior.ColorIndex = $HEADER_COLOR_INDEX
}

# ---------------------------------------------------------------------
# Single-column sheets: Company / Contact / Country.
# Each has header in col A/B (col A = the data header, col B =
# "FineshedAt"), data rows 2-7, and the finish timestamp written into
# row 4 of column B. We refresh that timestamp, then duplicate the
# A:B block into C:D with a brand new timestamp.
# ---------------------------------------------------------------------

function Write-TwoColTable($sheetName, $header, $values, $ts1, $ts2, $colWidth) {
    $ws = $wb.Worksheets.Item($sheetName)

    # refresh the timestamp of the existing run (B4), leave A:B data as is
    $ws.Cells.Item(4, 2).Value = $ts1

    # duplicate header (C1:D1)
    $ws.Cells.Item(1, 3).Value = $header
    Style-Header($ws.Cells.Item(1, 3))
    $ws.Cells.Item(1, 4).Value = "FineshedAt"
    Style-Header($ws.Cells.Item(1, 4))

    # duplicate data rows (C2:C7)
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($i + 2, 3).Value = $values[$i]
    }

    # duplicate finish timestamp (D4)
    $ws.Cells.Item(4, 4).Value = $ts2

    $ws.Columns.Item(3).ColumnWidth = $colWidth
    $ws.Columns.Item(4).ColumnWidth = $DATE_COL_WIDTH
}

$companyValues = @(
    "Alfreds Futterkiste",
    "Centro comercial Moctezuma",
    "Ernst Handel",
    "Island Trading",
    "Laughing Bacchus Winecellars",
    "Magazzini Alimentari Riuniti"
)
Write-TwoColTable "writeCompanyColumnIntoXcel" "Company" $companyValues `
    "1571855800605`nWed Oct 23 11:36:40 PDT 2019" `
    "1571855967414`nWed Oct 23 11:39:27 PDT 2019" `
    $WIDE_COL_WIDTH

$contactValues = @(
    "Maria Anders",
    "Francisco Chang",
    "Roland Mendel",
    "Helen Bennett",
    "Yoshi Tannamuri",
    "Giovanni Rovelli"
)
Write-TwoColTable "writeContactColumnIntoXcel" "Contact" $contactValues `
    "1571855802576`nWed Oct 23 11:36:42 PDT 2019" `
    "1571855969248`nWed Oct 23 11:39:29 PDT 2019" `
    $MEDIUM_COL_WIDTH

$countryValues = @(
    "Germany",
    "Mexico",
    "Austria",
    "UK",
    "Canada",
    "Italy"
)
Write-TwoColTable "writeCountryColumnIntoXcel" "Country" $countryValues `
    "1571855803803`nWed Oct 23 11:36:43 PDT 2019" `
    "1571855971281`nWed Oct 23 11:39:31 PDT 2019" `
    $NARROW_COL_WIDTH

# ---------------------------------------------------------------------
# Whole-table sheets: Company / Contact / Country / FineshedAt across
# A:D, duplicated into E:H with a fresh timestamp.
# ---------------------------------------------------------------------

function Write-FourColTable($sheetName, $ts1, $ts2) {
    $ws = $wb.Worksheets.Item($sheetName)

    # refresh the timestamp of the existing run (D4)
    $ws.Cells.Item(4, 4).Value = $ts1

    $headers = @("Company", "Contact", "Country", "FineshedAt")
    for ($c = 0; $c -lt 4; $c++) {
        $ws.Cells.Item(1, 5 + $c).Value = $headers[$c]
        Style-Header($ws.Cells.Item(1, 5 + $c))
    }

    for ($r = 2; $r -le 7; $r++) {
        $ws.Cells.Item($r, 5).Value = $ws.Cells.Item($r, 1).Value()
        $ws.Cells.Item($r, 6).Value = $ws.Cells.Item($r, 2).Value()
        $ws.Cells.Item($r, 7).Value = $ws.Cells.Item($r, 3).Value()
    }

    # duplicate finish timestamp (H4)
    $ws.Cells.Item(4, 8).Value = $ts2

    $ws.Columns.Item(5).ColumnWidth = $WIDE_COL_WIDTH
    $ws.Columns.Item(6).ColumnWidth = $MEDIUM_COL_WIDTH
    $ws.Columns.Item(7).ColumnWidth = $NARROW_COL_WIDTH
    $ws.Columns.Item(8).ColumnWidth = $DATE_COL_WIDTH
}

Write-FourColTable "writeWholeTableNestedFor" `
    "1571855810696`nWed Oct 23 11:36:50 PDT 2019" `
    "1571855997720`nWed Oct 23 11:39:57 PDT 2019"
Write-FourColTable "writeWholeTableSingleFor" `
    "1571855812034`nWed Oct 23 11:36:52 PDT 2019" `
    "1571856027433`nWed Oct 23 11:40:27 PDT 2019"

Write-Host "edit complete"
